$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Finish 4.6 through Data Structures" - log today's study session in the
# next open row of the time-log table (row 77: 2021-11-11).
$ws.Range("B77").Value = "JS101: Programming Foundations with JavaScript"
$ws.Range("C77").Value = 0.5
$ws.Range("D77").Value = "Lesson 4.6: Problem through Data Structures"

# Leave the selection where the user ended up next (matches the saved view).
$ws.Range("C78").Select()
